$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Cells.Item(60, 1)
$cell.Interior.ThemeColor = 5
$cell.Interior.TintAndShade = 0.4
Write-Host ("TAS=" + $cell.Interior.TintAndShade)
Write-Host ("Color=" + $cell.Interior.Color)
